$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in Test Steps (F3) and Test Data (G3) ---
$ws.Range("F3").Value = "1. Click on 'Register Here'`n2. Enter first name`n3. Enter last name`n4. Enter phone number`n5. Enter birthdate`n6. Enter email`n7. Enter password"
$ws.Range("G3").Value = "first name = Ben`nlast name = Dover`nphone number = 123456789`nbirthdate = 22-01-1993`nemail = customer@gmail.com`npassword = customer"

# --- Row 4: fill in Test Steps (F4), Test Data (G4), Expected Results (H4), Actual Results (I4), Pass/Fail (J4) ---
$ws.Range("F4").Value = "1. Click on 'Register Here'`n2. Enter first name`n3. Enter last name`n4. Enter phone number`n5. Enter birthdate`n6. Enter email`n7. Enter password"
$ws.Range("G4").Value = "first name = Ben`nlast name = Dover`nphone number = 12345678910`nbirthdate = 22-01-2500`nemail = customer`npassword = 123"
$ws.Range("H4").Value = "Register should not be successful and there should be specific error messages saying what is wrong"
$ws.Range("I4").Value = "As `nExpected"
$ws.Range("J4").Value = "Pass"

# --- Row 19: fill in Pre-Condition (E19), Test Steps (F19), Test Data (G19), Expected Results (H19), Actual Results (I19), Pass/Fail (J19) ---
$ws.Range("E19").Value = "user is registered`n and user information `nis stored in the `ndatabase"
$ws.Range("F19").Value = "1. Click on 'Edit Profile' found under 'My Profile'`n2. Enter first name`n3. Enter last name`n4. Enter phone number`n5. Enter birthdate`n6. Enter email`n7. Enter password"
$ws.Range("G19").Value = "first name = Saul`nlast name = T. Nutz`nphone number = 987654321`nbirthdate = 22-01-1994`nemail = customer1@gmail.com`npassword = customer1"
$ws.Range("H19").Value = "User info is changed"
$ws.Range("I19").Value = "As `nExpected"
$ws.Range("J19").Value = "Pass"

# --- Row 20: fill in Pre-Condition (E20), Test Steps (F20), Test Data (G20), Expected Results (H20), Actual Results (I20), Pass/Fail (J20) ---
$ws.Range("E20").Value = "user is registered`n and user information `nis stored in the `ndatabase"
$ws.Range("F20").Value = "1. Click on 'Edit Profile' found under 'My Profile'`n2. Enter first name`n3. Enter last name`n4. Enter phone number`n5. Enter birthdate`n6. Enter email`n7. Enter password"
$ws.Range("G20").Value = "(The same information as when registered or previously changed)`nfirst name = Ben`nlast name = Dover`nphone number = 123456789`nbirthdate = 22-01-1993`nemail = customer@gmail.com`npassword = customer"
$ws.Range("H20").Value = "User info is changed"
$ws.Range("I20").Value = "As `nExpected"
$ws.Range("J20").Value = "Pass"

# --- Apply wrap-text style (style index 1) to all newly-populated cells, matching source sheet formatting ---
$ws.Range("F3:G3").WrapText = $true
$ws.Range("F4:J4").WrapText = $true
$ws.Range("E19:J19").WrapText = $true
$ws.Range("E20:J20").WrapText = $true

# --- Row heights ---
$ws.Range("A3").RowHeight = 188.5
$ws.Range("A4").RowHeight = 188.5
$ws.Range("A19").RowHeight = 217.5
$ws.Range("A20").RowHeight = 217.5

# --- Column G width change ---
$ws.Range("G1").ColumnWidth = 27.1666666666667

# --- Selection change on the sheet (activeCell G20 -> H20) ---
$ws.Range("H20").Select()
